# Auto-generated edit script: updates numeric cell values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2880.25
$ws.Range("I53").Value = 453.875
$ws.Range("J53").Value = 7733
$ws.Range("K53").Value = 453.875
$ws.Range("L53").Value = 7733
$ws.Range("M53").Value = 183.125
$ws.Range("N53").Value = -9007
$ws.Range("H62").Value = 3029.1155
$ws.Range("I62").Value = 2064.6428
$ws.Range("K62").Value = 2064.6428
$ws.Range("M62").Value = -1440.6428
$ws.Range("H64").Value = 3207.9167
$ws.Range("I64").Value = 2999
$ws.Range("J64").Value = 3357.1428
$ws.Range("K64").Value = 2999
$ws.Range("L64").Value = 3357.1428
$ws.Range("M64").Value = -2751
$ws.Range("N64").Value = -3853.1428
$ws.Range("H65").Value = 3029.1155
$ws.Range("I65").Value = 2064.6428
$ws.Range("K65").Value = 10323.214
$ws.Range("M65").Value = -7203.214
$ws.Range("H67").Value = 3207.9167
$ws.Range("I67").Value = 2999
$ws.Range("J67").Value = 3357.1428
$ws.Range("K67").Value = 2999
$ws.Range("L67").Value = 3357.1428
$ws.Range("M67").Value = -2141
$ws.Range("N67").Value = -5073.1428
$ws.Range("H76").Value = 3706617
$ws.Range("I76").Value = 2914.2856
$ws.Range("J76").Value = 6947357
$ws.Range("K76").Value = 2914.2856
$ws.Range("L76").Value = 6947357
$ws.Range("M76").Value = -2599.2856
$ws.Range("N76").Value = -6947987
$ws.Range("H79").Value = 3706617
$ws.Range("I79").Value = 2914.2856
$ws.Range("J79").Value = 6947357
$ws.Range("K79").Value = 2914.2856
$ws.Range("L79").Value = 6947357
$ws.Range("M79").Value = -1822.2856
$ws.Range("N79").Value = -6949541
$ws.Range("H86").Value = 8641.429
$ws.Range("I86").Value = 1569.1428
$ws.Range("K86").Value = 1569.1428
$ws.Range("M86").Value = -446.1428000000001
$ws.Range("H89").Value = 8641.429
$ws.Range("I89").Value = 1569.1428
$ws.Range("K89").Value = 7845.714
$ws.Range("M89").Value = -2229.714
$ws.Range("H129").Value = 1438.2667
$ws.Range("J129").Value = 1474.069
$ws.Range("L129").Value = 4422.207
$ws.Range("N129").Value = -14422.207
$ws.Range("H132").Value = 2560.205
$ws.Range("I132").Value = 2364.4211
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 7093.263300000001
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -4563.263300000001
$ws.Range("N132").Value = -35060
$ws.Range("H137").Value = 104283.84
$ws.Range("I137").Value = 122368.695
$ws.Range("K137").Value = 367106.085
$ws.Range("M137").Value = -364556.085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 11000
$ws.Range("I23").Value = 25000
$ws.Range("K23").Value = 25000
$ws.Range("M23").Value = -24741
$ws.Range("H27").Value = 11201
$ws.Range("J27").Value = 11249.5
$ws.Range("L27").Value = 11249.5
$ws.Range("N27").Value = -11617.5
$ws.Range("H32").Value = 8103.927
$ws.Range("I32").Value = 5984.203
$ws.Range("J32").Value = 19354.77
$ws.Range("K32").Value = 5984.203
$ws.Range("L32").Value = 19354.77
$ws.Range("M32").Value = -5697.203
$ws.Range("N32").Value = -19928.77
$ws.Range("H63").Value = 3127646
$ws.Range("I63").Value = 2957.5
$ws.Range("J63").Value = 15626400
$ws.Range("K63").Value = 2957.5
$ws.Range("L63").Value = 15626400
$ws.Range("M63").Value = -2271.5
$ws.Range("N63").Value = -15627772
$ws.Range("H66").Value = 3127646
$ws.Range("I66").Value = 2957.5
$ws.Range("J66").Value = 15626400
$ws.Range("K66").Value = 14787.5
$ws.Range("L66").Value = 78132000
$ws.Range("M66").Value = -11355.5
$ws.Range("N66").Value = -78138864
$ws.Range("H97").Value = 495.18182
$ws.Range("I97").Value = 495.18182
$ws.Range("K97").Value = 495.18182
$ws.Range("M97").Value = 0.8181799999999839
$ws.Range("H102").Value = 1297.826
$ws.Range("I102").Value = 1192.5
$ws.Range("K102").Value = 1192.5
$ws.Range("M102").Value = 429.5
$ws.Range("H132").Value = 17774.363
$ws.Range("I132").Value = 2244.56
$ws.Range("K132").Value = 6733.68
$ws.Range("M132").Value = -4203.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3747.8857
$ws.Range("I134").Value = 3799.2942
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 11397.8826
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -8862.882599999999
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 144.32
$ws.Range("I22").Value = 144.875
$ws.Range("J22").Value = 143.33333
$ws.Range("K22").Value = 144.875
$ws.Range("L22").Value = 143.33333
$ws.Range("M22").Value = 205.125
$ws.Range("N22").Value = -843.3333299999999
$ws.Range("H58").Value = 64992.5
$ws.Range("I58").Value = 3228
$ws.Range("J58").Value = 126757
$ws.Range("K58").Value = 3228
$ws.Range("L58").Value = 126757
$ws.Range("M58").Value = -3025
$ws.Range("N58").Value = -127163
$ws.Range("H94").Value = 4281.769
$ws.Range("J94").Value = 6879.8
$ws.Range("L94").Value = 6879.8
$ws.Range("N94").Value = -7781.8
$ws.Range("H105").Value = 7353597
$ws.Range("I105").Value = 9615957
$ws.Range("J105").Value = 927.75
$ws.Range("K105").Value = 9615957
$ws.Range("L105").Value = 927.75
$ws.Range("M105").Value = -9614210
$ws.Range("N105").Value = -4421.75
$ws.Range("H132").Value = 1903.6744
$ws.Range("I132").Value = 1340.7878
$ws.Range("K132").Value = 4022.3634
$ws.Range("M132").Value = -1492.3634
$ws.Range("H134").Value = 1565.125
$ws.Range("I134").Value = 1164
$ws.Range("J134").Value = 2233.6667
$ws.Range("K134").Value = 3492
$ws.Range("L134").Value = 6701.000100000001
$ws.Range("M134").Value = -957
$ws.Range("N134").Value = -11771.0001
$ws.Range("H136").Value = 64992.5
$ws.Range("I136").Value = 3228
$ws.Range("J136").Value = 126757
$ws.Range("K136").Value = 9684
$ws.Range("L136").Value = 380271
$ws.Range("M136").Value = -7134
$ws.Range("N136").Value = -385371
$ws.Range("H141").Value = 27882
$ws.Range("J141").Value = 27882
$ws.Range("L141").Value = 27882
$ws.Range("N141").Value = -38242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 43
$ws.Range("J33").Value = 59
$ws.Range("L33").Value = 354
$ws.Range("N33").Value = -920
$ws.Range("H51").Value = 3363.5
$ws.Range("I51").Value = 1800
$ws.Range("K51").Value = 5400
$ws.Range("M51").Value = -4940
$ws.Range("H74").Value = 9900
$ws.Range("J74").Value = 9900
$ws.Range("L74").Value = 29700
$ws.Range("N74").Value = -31822
$ws.Range("H75").Value = 2158.5715
$ws.Range("J75").Value = 2182.8333
$ws.Range("L75").Value = 6548.499899999999
$ws.Range("N75").Value = -8544.499899999999
$ws.Range("H77").Value = 9900
$ws.Range("J77").Value = 9900
$ws.Range("L77").Value = 89100
$ws.Range("N77").Value = -99708
$ws.Range("H78").Value = 2158.5715
$ws.Range("J78").Value = 2182.8333
$ws.Range("L78").Value = 19645.4997
$ws.Range("N78").Value = -29629.4997
$ws.Range("H82").Value = 10000
$ws.Range("J82").Value = 10000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 10000
$ws.Range("J85").Value = 10000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32808
$ws.Range("H100").Value = 3172.8
$ws.Range("J100").Value = 3172.8
$ws.Range("L100").Value = 9518.400000000001
$ws.Range("N100").Value = -11140.4
$ws.Range("H117").Value = 1497.6666
$ws.Range("I117").Value = 895.8
$ws.Range("J117").Value = 2250
$ws.Range("K117").Value = 2687.4
$ws.Range("L117").Value = 6750
$ws.Range("M117").Value = 754.6000000000004
$ws.Range("N117").Value = -13634
$ws.Range("H122").Value = 815.0476
$ws.Range("I122").Value = 316.33334
$ws.Range("J122").Value = 1014.5333
$ws.Range("K122").Value = 2847.00006
$ws.Range("L122").Value = 9130.7997
$ws.Range("M122").Value = -397.0000600000003
$ws.Range("N122").Value = -14030.7997
$ws.Range("H131").Value = 673.9400000000001
$ws.Range("J131").Value = 751.8570999999999
$ws.Range("L131").Value = 2255.5713
$ws.Range("N131").Value = -12335.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1093218
$ws.Range("I122").Value = 1512301.4
$ws.Range("K122").Value = 4536904.199999999
$ws.Range("M122").Value = -4534454.199999999
$ws.Range("H136").Value = 2878.8667
$ws.Range("I136").Value = 2488.4
$ws.Range("J136").Value = 3659.8
$ws.Range("K136").Value = 7465.200000000001
$ws.Range("L136").Value = 10979.4
$ws.Range("M136").Value = -4915.200000000001
$ws.Range("N136").Value = -16079.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1222.8
$ws.Range("I132").Value = 486.91666
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 1460.74998
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = 1069.25002
$ws.Range("N132").Value = -17559.0005
$ws.Range("H136").Value = 29496140
$ws.Range("I136").Value = 36867676
$ws.Range("K136").Value = 110603028
$ws.Range("M136").Value = -110600478
